# Update "想去人数" (F2) and "最低票价" (G2) values on the sheets that
# contain this event row: "展览" and "全部类型".
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 533
    $ws.Range("G2").Value = 55
}
